# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G ("K") for rows 2-12
$newK = @{
    2  = 1
    3  = 0
    4  = 2
    5  = 1
    6  = 0
    7  = 0
    8  = 0
    9  = 0
    10 = 0
    11 = 2
    12 = 0
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
